# Apply the "removed comments from DoorController.rebeca" status update
# to the StatusTableOverview sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("StatusTableOverview")

# Update the "Code Edits" cell for the "Train Door Controllers" row (row 10)
# to reflect that comments were removed, and mark it with the "Good" cell style.
$ws.Range("E10").Value = "removed comments"
$ws.Range("E10").Style = "Good"

# Move the active selection, matching the recorded cursor position after the edit.
$ws.Range("F11").Select()
